$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# EnableScheduler row: mark status as Automated (finished testing), bump automated count
$ws.Range("D3").Value = "Automated"
$ws.Range("B3").Value = 5

# RequiredFields row: rename from "Suited to Manual" to "RequiredFields",
# finished testing (status already "Automated"), remove the old "Visuals" note
$ws.Range("A5").Value = "RequiredFields"
$ws.Range("E5").Clear()

# Update the active selection to D4
$ws.Range("D4").Select()

$wb.Save()
